$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 152.28572
$ws.Range("I5").Value = 169.2
$ws.Range("K5").Value = 169.2
$ws.Range("M5").Value = -54.19999999999999
$ws.Range("H22").Value = 180
$ws.Range("J22").Value = 180
$ws.Range("L22").Value = 540
$ws.Range("N22").Value = -884
$ws.Range("H30").Value = 976
$ws.Range("I30").Value = 976
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2928
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -2827
$ws.Range("N30").ClearContents()
$ws.Range("H46").Value = 3900
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3900
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 11700
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -11938
$ws.Range("H60").Value = 3900
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 3900
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 11700
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -12668
$ws.Range("H100").Value = 1976.421
$ws.Range("I100").Value = 1863.4667
$ws.Range("J100").Value = 2400
$ws.Range("K100").Value = 1863.4667
$ws.Range("L100").Value = 2400
$ws.Range("M100").Value = -1322.4667
$ws.Range("N100").Value = -3482
$ws.Range("H112").Value = 2969.2415
$ws.Range("J112").Value = 3007.4285
$ws.Range("L112").Value = 9022.2855
$ws.Range("N112").Value = -11238.2855
$ws.Range("H125").Value = 839.8823
$ws.Range("I125").Value = 591.4286
$ws.Range("J125").Value = 1013.8
$ws.Range("K125").Value = 5322.8574
$ws.Range("L125").Value = 9124.199999999999
$ws.Range("M125").Value = -2862.8574
$ws.Range("N125").Value = -14044.2
$ws.Range("H135").Value = 936.0789
$ws.Range("I135").Value = 936.0789
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 8424.7101
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -5889.7101
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 2531.0952
$ws.Range("I138").Value = 2152.6191
$ws.Range("J138").Value = 2909.5715
$ws.Range("K138").Value = 6457.8573
$ws.Range("L138").Value = 8728.7145
$ws.Range("M138").Value = -1317.8573
$ws.Range("N138").Value = -19008.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 381041.97
$ws.Range("I32").Value = 454612.62
$ws.Range("K32").Value = 454612.62
$ws.Range("M32").Value = -454325.62
$ws.Range("H61").Value = 2028.5962
$ws.Range("I61").Value = 1232.0769
$ws.Range("J61").Value = 2825.1155
$ws.Range("K61").Value = 1232.0769
$ws.Range("L61").Value = 2825.1155
$ws.Range("M61").Value = -1020.0769
$ws.Range("N61").Value = -3249.1155
$ws.Range("H74").Value = 1500.9375
$ws.Range("I74").Value = 1277.75
$ws.Range("J74").Value = 1575.3334
$ws.Range("K74").Value = 1277.75
$ws.Range("L74").Value = 1575.3334
$ws.Range("M74").Value = -403.75
$ws.Range("N74").Value = -3323.3334
$ws.Range("H77").Value = 1500.9375
$ws.Range("I77").Value = 1277.75
$ws.Range("J77").Value = 1575.3334
$ws.Range("K77").Value = 6388.75
$ws.Range("L77").Value = 7876.666999999999
$ws.Range("M77").Value = -2020.75
$ws.Range("N77").Value = -16612.667
$ws.Range("H132").Value = 5716.7354
$ws.Range("I132").Value = 6908.143
$ws.Range("J132").Value = 3792.1538
$ws.Range("K132").Value = 20724.429
$ws.Range("L132").Value = 11376.4614
$ws.Range("M132").Value = -18194.429
$ws.Range("N132").Value = -16436.4614
$ws.Range("H136").Value = 2028.5962
$ws.Range("I136").Value = 1232.0769
$ws.Range("J136").Value = 2825.1155
$ws.Range("K136").Value = 3696.2307
$ws.Range("L136").Value = 8475.3465
$ws.Range("M136").Value = -1146.2307
$ws.Range("N136").Value = -13575.3465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H92").Value = 70050.75
$ws.Range("J92").Value = 70050.75
$ws.Range("L92").Value = 70050.75
$ws.Range("N92").Value = -75042.75
$ws.Range("H99").Value = 1078.4445
$ws.Range("I99").Value = 981.7619
$ws.Range("J99").Value = 1416.8334
$ws.Range("K99").Value = 981.7619
$ws.Range("L99").Value = 1416.8334
$ws.Range("M99").Value = 516.2381
$ws.Range("N99").Value = -4412.8334
$ws.Range("H100").Value = 100000
$ws.Range("J100").Value = 100000
$ws.Range("L100").Value = 100000
$ws.Range("N100").Value = -102164
$ws.Range("H103").Value = 58000
$ws.Range("J103").Value = 58000
$ws.Range("L103").Value = 58000
$ws.Range("N103").Value = -60344
$ws.Range("H112").Value = 100469
$ws.Range("J112").Value = 100469
$ws.Range("L112").Value = 100469
$ws.Range("N112").Value = -103423

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 467.45
$ws.Range("I22").Value = 147
$ws.Range("J22").Value = 948.125
$ws.Range("K22").Value = 147
$ws.Range("L22").Value = 948.125
$ws.Range("M22").Value = 203
$ws.Range("N22").Value = -1648.125
$ws.Range("H31").Value = 4530.4346
$ws.Range("I31").Value = 1170.6316
$ws.Range("J31").Value = 6894.7407
$ws.Range("K31").Value = 1170.6316
$ws.Range("L31").Value = 6894.7407
$ws.Range("M31").Value = -875.6315999999999
$ws.Range("N31").Value = -7484.7407
$ws.Range("H34").Value = 4530.4346
$ws.Range("I34").Value = 1170.6316
$ws.Range("J34").Value = 6894.7407
$ws.Range("K34").Value = 1170.6316
$ws.Range("L34").Value = 6894.7407
$ws.Range("M34").Value = -968.6315999999999
$ws.Range("N34").Value = -7298.7407
$ws.Range("H94").Value = 1207.7858
$ws.Range("I94").Value = 661.5
$ws.Range("J94").Value = 1426.3
$ws.Range("K94").Value = 661.5
$ws.Range("L94").Value = 1426.3
$ws.Range("M94").Value = -210.5
$ws.Range("N94").Value = -2328.3
$ws.Range("H132").Value = 3970297
$ws.Range("I132").Value = 1572.5
$ws.Range("J132").Value = 8335893.5
$ws.Range("K132").Value = 4717.5
$ws.Range("L132").Value = 25007680.5
$ws.Range("M132").Value = -2187.5
$ws.Range("N132").Value = -25012740.5
$ws.Range("H134").Value = 1170.4
$ws.Range("I134").Value = 784
$ws.Range("J134").Value = 1750
$ws.Range("K134").Value = 2352
$ws.Range("L134").Value = 5250
$ws.Range("M134").Value = 183
$ws.Range("N134").Value = -10320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1286.5408
$ws.Range("I68").Value = 742.4
$ws.Range("J68").Value = 1472.8904
$ws.Range("K68").Value = 2227.2
$ws.Range("L68").Value = 4418.6712
$ws.Range("M68").Value = -1416.2
$ws.Range("N68").Value = -6040.6712
$ws.Range("H71").Value = 1286.5408
$ws.Range("I71").Value = 742.4
$ws.Range("J71").Value = 1472.8904
$ws.Range("K71").Value = 6681.599999999999
$ws.Range("L71").Value = 13256.0136
$ws.Range("M71").Value = -2625.599999999999
$ws.Range("N71").Value = -21368.0136
$ws.Range("H131").Value = 1151.2307
$ws.Range("I131").Value = 994.75
$ws.Range("J131").Value = 1191.6129
$ws.Range("K131").Value = 2984.25
$ws.Range("L131").Value = 3574.8387
$ws.Range("M131").Value = 2055.75
$ws.Range("N131").Value = -13654.8387

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 22826494
$ws.Range("I80").Value = 38040092
$ws.Range("J80").Value = 6099
$ws.Range("K80").Value = 38040092
$ws.Range("L80").Value = 6099
$ws.Range("M80").Value = -38039094
$ws.Range("N80").Value = -8095
$ws.Range("H83").Value = 22826494
$ws.Range("I83").Value = 38040092
$ws.Range("J83").Value = 6099
$ws.Range("K83").Value = 190200460
$ws.Range("L83").Value = 30495
$ws.Range("M83").Value = -190195468
$ws.Range("N83").Value = -40479
$ws.Range("H132").Value = 1981.6666
$ws.Range("I132").Value = 1408.7059
$ws.Range("J132").Value = 3373.1428
$ws.Range("K132").Value = 4226.1177
$ws.Range("L132").Value = 10119.4284
$ws.Range("M132").Value = -1696.1177
$ws.Range("N132").Value = -15179.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 50002724
$ws.Range("I7").Value = 83335860
$ws.Range("J7").Value = 3026.25
$ws.Range("K7").Value = 83335860
$ws.Range("L7").Value = 3026.25
$ws.Range("M7").Value = -83335748
$ws.Range("N7").Value = -3250.25
$ws.Range("H68").Value = 1390.1305
$ws.Range("I68").Value = 1267.0625
$ws.Range("J68").Value = 1671.4286
$ws.Range("K68").Value = 1267.0625
$ws.Range("L68").Value = 1671.4286
$ws.Range("M68").Value = -518.0625
$ws.Range("N68").Value = -3169.4286
$ws.Range("H71").Value = 1390.1305
$ws.Range("I71").Value = 1267.0625
$ws.Range("J71").Value = 1671.4286
$ws.Range("K71").Value = 6335.3125
$ws.Range("L71").Value = 8357.143
$ws.Range("M71").Value = -2591.3125
$ws.Range("N71").Value = -15845.143
$ws.Range("H122").Value = 7500.5713
$ws.Range("I122").Value = 6502
$ws.Range("K122").Value = 19506
$ws.Range("M122").Value = -17056
$ws.Range("H126").Value = 50002724
$ws.Range("I126").Value = 83335860
$ws.Range("J126").Value = 3026.25
$ws.Range("K126").Value = 250007580
$ws.Range("L126").Value = 9078.75
$ws.Range("M126").Value = -250005110
$ws.Range("N126").Value = -14018.75
$ws.Range("H132").Value = 4150.079
$ws.Range("I132").Value = 4015.682
$ws.Range("J132").Value = 4334.875
$ws.Range("K132").Value = 12047.046
$ws.Range("L132").Value = 13004.625
$ws.Range("M132").Value = -9517.045999999998
$ws.Range("N132").Value = -18064.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 815.1539
$ws.Range("I100").Value = 633.8570999999999
$ws.Range("K100").Value = 1267.7142
$ws.Range("M100").Value = -726.7141999999999
$ws.Range("H132").Value = 3970361
$ws.Range("I132").Value = 1939.9584
$ws.Range("J132").Value = 9261589
$ws.Range("K132").Value = 5819.8752
$ws.Range("L132").Value = 27784767
$ws.Range("M132").Value = -3289.8752
$ws.Range("N132").Value = -27789827
$ws.Range("H136").Value = 2692.98
$ws.Range("I136").Value = 2762.8215
$ws.Range("J136").Value = 2604.0908
$ws.Range("K136").Value = 8288.4645
$ws.Range("L136").Value = 7812.2724
$ws.Range("M136").Value = -5738.4645
$ws.Range("N136").Value = -12912.2724
